$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Add new row 3 data, mirroring the formula pattern already present in row 2.
$ws.Range("A3").Value = 1.05
$ws.Range("B3").Value = 11.9
$ws.Range("C3").Formula = "=B3*A3"
$ws.Range("D3").Value = 1.054
$ws.Range("E3").Formula = "=C3/D3"
$ws.Range("F3").Formula = "=E3*D3"
$ws.Range("G3").Formula = "=B3-E3"

$ws.Activate()
$ws.Range("E3").Select()
